# Auto-generated edit script applying the crypto price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.662.15"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").Value = "1.962.80"
$ws.Range("E3").Value = "  +2.07%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.92"
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  +2.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.50"
$ws.Range("E7").Value = "  +9.44%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.378"
$ws.Range("E9").Value = "  +7.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0792"
$ws.Range("E10").Value = "  -3.02%  "
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.20"
$ws.Range("E12").Value = "  +9.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.841"
$ws.Range("E13").Value = "  +7.61%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.251.79"
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.71"
$ws.Range("E15").Value = "  +6.38%  "
$ws.Range("E16").Value = "  +5.36%  "
$ws.Range("D17").Value = "1.949.15"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").Value = "36.585.34"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.91"
$ws.Range("E19").Value = "  +2.96%  "
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "230.39"
$ws.Range("E21").Value = "  +2.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.10"
$ws.Range("E22").Value = "  +5.34%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.47"
$ws.Range("E24").Value = "  +7.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("E25").Value = "  +5.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.145"
$ws.Range("E26").Value = "  +12.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.24"
$ws.Range("E27").Value = "  +3.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.85"
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.49"
$ws.Range("E29").Value = "  +3.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.22"
$ws.Range("E30").Value = "  +13.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("E31").Value = "  +3.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.77"
$ws.Range("E32").Value = "  +8.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0618"
$ws.Range("E33").Value = "  +1.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.43"
$ws.Range("E34").Value = "  +9.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.58"
$ws.Range("E35").Value = "  +24.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.28"
$ws.Range("E36").Value = "  +9.84%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.76"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.59"
$ws.Range("E39").Value = "  -6.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0983"
$ws.Range("E40").Value = "  +3.24%  "
$ws.Range("E41").Value = "  +1.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.17"
$ws.Range("E42").Value = "  +5.13%  "
$ws.Range("E43").Value = "  +3.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.16"
$ws.Range("E44").Value = "  +6.48%  "
$ws.Range("D45").Value = "1.370.02"
$ws.Range("E45").Value = "  +4.41%  "
$ws.Range("E46").Value = "  +4.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.52"
$ws.Range("E47").Value = "  +6.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.17"
$ws.Range("E48").Value = "  +3.38%  "
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.42"
$ws.Range("E50").Value = "  +5.45%  "
$ws.Range("D51").Value = "2.142.05"
$ws.Range("E51").Value = "  +1.90%  "
